$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 172-174), mirroring the existing schema
$rows = @(
    @{ A=3; B="Femacal de La Calera"; C="Coquimbo"; D=44911; E=5; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103003; J="Damasco"; K="Castle Brite"; L="Especial"; M=65; N=15000; O=15000; P=15000; Q="$/caja 15 kilos"; R="Provincia de San Felipe de Aconcagua"; S=1000; T=15 },
    @{ A=3; B="Femacal de La Calera"; C="Coquimbo"; D=44911; E=5; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103003; J="Damasco"; K="Castle Brite"; L="Primera"; M=64; N=13000; O=13000; P=13000; Q="$/caja 15 kilos"; R="Provincia de San Felipe de Aconcagua"; S=867; T=15 },
    @{ A=3; B="Femacal de La Calera"; C="Coquimbo"; D=44911; E=5; F="Fruta"; G=100103; H="Frutos de hueso (carozo)"; I=100103003; J="Damasco"; K="Castle Brite"; L="Segunda"; M=60; N=11000; O=11000; P=11000; Q="$/caja 15 kilos"; R="Provincia de San Felipe de Aconcagua"; S=733; T=15 }
)

$startRow = 172

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $data.D
    $dCell.NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
}
